$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-09-02 12:52:18"

$wsZhCn.Range("H4").Value = "2016-09-02 12:52:13"
$wsZhCn.Range("K4").Value = "2016-09-02 12:52:41"

$wsDeDe.Range("H4").Value = "2016-09-02 12:52:18"
$wsDeDe.Range("K4").Value = "2016-09-02 12:52:48"
